$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Numeric-looking Price values are prefixed with a leading apostrophe so Excel
# stores them as text (matching the original inline-string format) instead of
# auto-converting them to numbers and losing formatting (e.g. trailing zeros).

$ws.Range("D2").Value = "27.725.32"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.859.71"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +0.88%  "
$ws.Range("D5").Value = "'332.94"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").Value = "'0.4698"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").Value = "'0.3898"
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").Value = "'46.63"
$ws.Range("D10").Value = "'0.07991"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("D12").Value = "'21.59"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").Value = "1.873.15"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").Value = "'5.999"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "'7.139"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "'1.014"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "'0.06697"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'0.00001040"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "'16.90"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "27.735.91"
$ws.Range("D23").Value = "'5.461"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "'2.315"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").Value = "2.095.09"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").Value = "'158.52"
$ws.Range("D28").Value = "'19.67"
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("D29").Value = "'2.091"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").Value = "'5.399"
$ws.Range("E30").Value = "  -2.93%  "
$ws.Range("D31").Value = "'120.81"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").Value = "'0.9671"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").Value = "'0.09436"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D34").Value = "'3.637"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("D35").Value = "'5.307"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").Value = "'1.343"
$ws.Range("E36").Value = "  -7.20%  "
$ws.Range("D37").Value = "'0.06035"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "'0.02213"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").Value = "'1.207"
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").Value = "'8.144"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").Value = "'0.5924"
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D44").Value = "'10.21"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'1.256"
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").Value = "'0.5615"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").Value = "'12.01"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").Value = "'1.917"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").Value = "'3.301"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").Value = "'0.06771"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").Value = "'112.76"
$ws.Range("E51").Value = "  -1.78%  "
